$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old column F ("ElementName3" header, last column of the old layout)
$ws.Range("F1").ClearContents()

# Fill in the new data for columns B:E, rows 2-9 (previously only column A was populated)
$ws.Range("B2").Value = "pede. Suspendisse dui."
$ws.Range("C2").Value = "EL Haddad"
$ws.Range("D2").Value = "Nullam feugiat placerat"
$ws.Range("E2").Value = "varius et, euismod"

$ws.Range("B3").Value = "a nunc. In"
$ws.Range("C3").Value = "Badir"
$ws.Range("D3").Value = "sodales nisi magna"
$ws.Range("E3").Value = "elementum sem, vitae"

$ws.Range("B4").Value = "amet metus. Aliquam"
$ws.Range("C4").Value = "Ezzine"
$ws.Range("D4").Value = "Cras vulputate velit"
$ws.Range("E4").Value = "scelerisque neque sed"

$ws.Range("B5").Value = "quam vel sapien"
$ws.Range("C5").Value = "El Alami Hassoun"
$ws.Range("D5").Value = "Nunc mauris elit,"
$ws.Range("E5").Value = "libero et tristique"

$ws.Range("B6").Value = "feugiat nec, diam."
$ws.Range("C6").Value = "Lazaar"
$ws.Range("D6").Value = "pellentesque. Sed dictum."
$ws.Range("E6").Value = "ridiculus mus. Proin"

$ws.Range("B7").Value = "nonummy. Fusce fermentum"
$ws.Range("C7").Value = "El Haddad"
$ws.Range("D7").Value = "neque pellentesque massa"
$ws.Range("E7").Value = "Mauris eu turpis."

$ws.Range("B8").Value = "a, arcu. Sed"
$ws.Range("C8").Value = "EL Haddad"
$ws.Range("D8").Value = "sit amet risus."
$ws.Range("E8").Value = "Nulla facilisi. Sed"

$ws.Range("B9").Value = "Suspendisse eleifend. Cras"
$ws.Range("C9").Value = "El Alami Hassoun"
$ws.Range("D9").Value = "velit dui, semper"
$ws.Range("E9").Value = "ligula elit, pretium"

# B2 came in with an explicit black font color (e.g. pasted from another source)
$ws.Range("B2").Font.Color = 0

# Resize the columns to fit their new, wider content (AutoFit-style behaviour).
# ColumnWidth is expressed in characters, excluding the ~5px/6 cell padding that
# Excel adds on top when it stores the column width in the sheet XML, so that
# padding is subtracted here to land on the desired rendered width.
$ws.Columns("A:A").ColumnWidth = 9.28515625 - 0.8333333333333334
$ws.Columns("B:B").ColumnWidth = 27 - 0.8333333333333334
$ws.Columns("C:C").ColumnWidth = 16.140625 - 0.8333333333333334
$ws.Columns("D:D").ColumnWidth = 25.140625 - 0.8333333333333334
$ws.Columns("E:E").ColumnWidth = 21.28515625 - 0.8333333333333334

# Move the selection to where it ended up after the edits
$ws.Range("H11").Select() | Out-Null
